$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 2 and 3
$wsOverview.Range("G2").Value = "2016-08-17 16:16:06"
$wsOverview.Range("G3").Value = "2016-08-17 16:16:06"

# zh-cn sheet: "Priority" column (E) ht -> mt, rows 2 and 3
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: "Correspond Handoff Datetime" column (H), rows 2 and 3
# (both rows shared the same underlying string in the original workbook)
$wsZhCn.Range("H2").Value = "2016-08-17 16:15:56"
$wsZhCn.Range("H3").Value = "2016-08-17 16:15:56"

# zh-cn sheet: "Correspond Handback DateTime" column (K), rows 2 and 3
$wsZhCn.Range("K2").Value = "2016-08-17 16:16:29"
$wsZhCn.Range("K3").Value = "2016-08-17 16:16:29"

# de-de sheet: "Priority" column (E) ht -> mt, rows 2 and 3
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# de-de sheet: "Correspond Handoff Datetime" column (H), rows 2 and 3
$wsDeDe.Range("H2").Value = "2016-08-17 16:16:06"
$wsDeDe.Range("H3").Value = "2016-08-17 16:16:06"

# de-de sheet: "Correspond Handback DateTime" column (K), rows 2 and 3
# (both rows shared the same underlying string in the original workbook)
$wsDeDe.Range("K2").Value = "2016-08-17 16:16:37"
$wsDeDe.Range("K3").Value = "2016-08-17 16:16:37"
